# Weekly price update: insert two new daily records into the Zanahoria
# (carrot) price series, each insertion pushing the rows below it down
# by one (matching the workbook's newest-date-at-top-of-block ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insertion #1: new record becomes row 146 -----------------------------
$ws.Rows(146).Insert()

$ws.Range("A146").Value = 7
$ws.Range("B146").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C146").Value = "Ñuble"
$ws.Range("D146").Value = 44567
$ws.Range("E146").Value = 16
$ws.Range("F146").Value = 100114013
$ws.Range("G146").Value = "Zanahoria"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 300
$ws.Range("K146").Value = 6500
$ws.Range("L146").Value = 7000
$ws.Range("M146").Value = 6750
$ws.Range("N146").Value = "`$/saco 20 kilos"
$ws.Range("O146").Value = "Provincia de Diguillín"
$ws.Range("P146").Value = 338
$ws.Range("Q146").Value = 20
$ws.Range("R146").Value = "Hortaliza"

# --- Insertion #2: new record becomes row 230 ------------------------------
$ws.Rows(230).Insert()

$ws.Range("A230").Value = 7
$ws.Range("B230").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C230").Value = "Ñuble"
$ws.Range("D230").Value = 44568
$ws.Range("E230").Value = 16
$ws.Range("F230").Value = 100114013
$ws.Range("G230").Value = "Zanahoria"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 100
$ws.Range("K230").Value = 6500
$ws.Range("L230").Value = 7000
$ws.Range("M230").Value = 6750
$ws.Range("N230").Value = "`$/saco 20 kilos"
$ws.Range("O230").Value = "Provincia de Diguillín"
$ws.Range("P230").Value = 338
$ws.Range("Q230").Value = 20
$ws.Range("R230").Value = "Hortaliza"
